# Applies the "Updated symbol list" price/volume/coin refresh described by
# the commit diff: columns B (Coin) and C (Link) get their text swapped for
# three newly-inserted/ranked coins (rows 8-17 shift by one rank), and
# columns D (Price) / E (Volume(1h)) get refreshed numeric/percentage
# readings across the affected rows.
#
# D/E values in this sheet are stored as plain text (t="inlineStr" in the
# original workbook, e.g. "306.18" / "0.92%"), not numbers - so we write
# them through Set-TextCell, which forces a leading apostrophe. That keeps
# Excel's "looks like a number/percentage" auto-detection from silently
# converting them to numeric cells (which would also pick up a numeric/
# percentage NumberFormat and no longer match the source text values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$Address,
        [string]$Text
    )
    # Leading "'" = Excel's classic "force text" entry marker; it is stored
    # as part of the quote-prefix flag, not as literal text in the cell.
    $ws.Range($Address).Value = "'" + $Text
}

Set-TextCell "D2" '306.22'
Set-TextCell "E2" '0.99%'

Set-TextCell "D3" '36.26'
Set-TextCell "E3" '-1.30%'

Set-TextCell "D4" '5.057'
Set-TextCell "E4" '1.34%'

Set-TextCell "D5" '0.07938'
Set-TextCell "E5" '2.77%'

Set-TextCell "D6" '2.234'
Set-TextCell "E6" '8.09%'

Set-TextCell "D7" '8.006'
Set-TextCell "E7" '0.90%'

$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell "D8" '0.9270'
Set-TextCell "E8" '1.35%'

$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextCell "D9" '0.09837'
Set-TextCell "E9" '2.96%'

$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell "D10" '0.1885'
Set-TextCell "E10" '1.81%'

$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextCell "D11" '0.09147'
Set-TextCell "E11" '7.01%'

$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextCell "D12" '0.03716'
Set-TextCell "E12" '4.21%'

$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell "D13" '0.09928'
Set-TextCell "E13" '-0.53%'

$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell "D14" '0.001442'
Set-TextCell "E14" '-2.04%'

$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextCell "D15" '0.005615'
Set-TextCell "E15" '-1.72%'

$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell "D16" '3.461'
Set-TextCell "E16" '0.02%'

$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextCell "D17" '4.141'
Set-TextCell "E17" '2.94%'

Set-TextCell "D18" '2.633'
Set-TextCell "E18" '18.64%'

Set-TextCell "E19" '-0.05%'

Set-TextCell "E20" '-0.93%'

Set-TextCell "D21" '5.099'
Set-TextCell "E21" '3.77%'

Set-TextCell "D22" '0.2248'
Set-TextCell "E22" '1.88%'

Set-TextCell "D23" '0.04547'

Set-TextCell "D24" '0.001239'
Set-TextCell "E24" '0.28%'

Set-TextCell "D25" '0.004783'
Set-TextCell "E25" '-6.18%'

Set-TextCell "D26" '0.0001300'
Set-TextCell "E26" '-7.10%'

Set-TextCell "E27" '73.83%'

Set-TextCell "D39" '0.01913'
Set-TextCell "E39" '9.34%'

Set-TextCell "D40" '0.04930'
Set-TextCell "E40" '6.78%'

Set-TextCell "D41" '0.007843'
Set-TextCell "E41" '2.21%'

Set-TextCell "D42" '0.1397'
Set-TextCell "E42" '0.31%'

Set-TextCell "D43" '0.007805'
Set-TextCell "E43" '0.87%'

Set-TextCell "D44" '0.002240'
Set-TextCell "E44" '3.66%'

Set-TextCell "D45" '0.01142'
Set-TextCell "E45" '10.07%'

Set-TextCell "D46" '0.00006274'
Set-TextCell "E46" '-1.21%'

Set-TextCell "D47" '0.00000000750'
Set-TextCell "E47" '-0.32%'

Set-TextCell "D48" '51.76'
Set-TextCell "E48" '49.88%'

Set-TextCell "D49" '0.001800'
Set-TextCell "E49" '-10.28%'

Set-TextCell "D50" '0.00002100'
Set-TextCell "E50" '-0.32%'

Set-TextCell "D51" '0.0002000'
Set-TextCell "E51" '-0.32%'
